$wb = $excel.ActiveWorkbook

# --- CustomerDetails sheet: replace the stray "sss" cell (F5) with a full
# --- second data row (A2:F2) of real customer details ---
$wsCustomer = $wb.Worksheets.Item("CustomerDetails")
$wsCustomer.Activate() | Out-Null

$wsCustomer.Range("F5").ClearContents() | Out-Null

$wsCustomer.Range("A2").Value = "Kasim"
$wsCustomer.Range("B2").Value = "mdkasim111@hotmail.com"
$wsCustomer.Range("D2").Value = "Jakarta"
$wsCustomer.Range("C2").Value = "'081808466410"
$wsCustomer.Range("E2").Value = "MidPlaza 2, 4th Floor Jl.Jend.Sudirman Kav.10-11"
$wsCustomer.Range("F2").Value = "'10220"

# Give the new columns an explicit (best-fit style) width, matching the
# columns that now hold longer text (Email, Phone, Address).
$wsCustomer.Columns.Item(2).AutoFit() | Out-Null
$wsCustomer.Columns.Item(3).AutoFit() | Out-Null
$wsCustomer.Columns.Item(5).AutoFit() | Out-Null

$wsCustomer.Range("E9").Select() | Out-Null

# --- FailedCCPayment sheet keeps its data (values unaffected aside from
# --- the shared-string renumbering caused by the edit above); only the
# --- user's current selection moved. Re-select it last so it stays the
# --- active sheet/tab, matching the original workbook state. ---
$wsFailed = $wb.Worksheets.Item("FailedCCPayment")
$wsFailed.Activate() | Out-Null
$wsFailed.Range("H28").Select() | Out-Null
